$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("27.429.36", "  -0.74%  ")
    3 = @("1.821.67", "  -2.45%  ")
    4 = @("1.003", "  -1.13%  ")
    5 = @("332.65", "  -0.85%  ")
    6 = @("1.003", "  -1.05%  ")
    7 = @("0.4581", "  -1.80%  ")
    8 = @("0.3811", "  -2.95%  ")
    9 = @("45.87", "  -0.31%  ")
    10 = @("0.07842", "  -1.84%  ")
    11 = @("0.9563", "  -4.55%  ")
    12 = @("21.00", "  -3.42%  ")
    13 = @("1.830.63", "  -1.97%  ")
    14 = @("5.831", "  -2.11%  ")
    15 = @("7.068", "  -2.21%  ")
    16 = @("1.004", "  -1.07%  ")
    17 = @("89.18", "  +0.67%  ")
    18 = @("0.06573", "  -2.06%  ")
    19 = @("0.00001018", "  -2.52%  ")
    20 = @("17.05", "  -1.11%  ")
    21 = @("1.002", "  -1.17%  ")
    22 = @("27.410.04", "  -0.83%  ")
    23 = @("5.285", "  -3.27%  ")
    24 = @("10.81", "  -1.20%  ")
    25 = @("2.255", "  -2.49%  ")
    26 = @("2.049.27", "  -2.07%  ")
    27 = @("157.18", "  -1.21%  ")
    28 = @("19.30", "  -1.86%  ")
    29 = @("2.039", "  -4.60%  ")
    30 = @("5.257", "  -3.37%  ")
    31 = @("117.60", "  -3.36%  ")
    32 = @("0.09313", "  -1.56%  ")
    33 = @("0.9297", "  -4.87%  ")
    34 = @("3.566", "  -1.58%  ")
    35 = @("5.214", "  -1.87%  ")
    36 = @("1.314", "  -2.00%  ")
    37 = @("0.05918", "  -1.89%  ")
    38 = @("0.02186", "  -2.26%  ")
    39 = @("8.080", "  -3.21%  ")
    40 = @("1.002", "  -1.12%  ")
    41 = @("1.141", "  -4.74%  ")
    42 = @("0.5739", "  -3.61%  ")
    43 = @("0.1814", "  -3.23%  ")
    44 = @("9.911", "  -4.06%  ")
    45 = @("1.274", "  +1.86%  ")
    46 = @("0.5386", "  -4.23%  ")
    47 = @("11.76", "  -3.17%  ")
    48 = @("1.871", "  -2.86%  ")
    49 = @("110.28", "  -1.43%  ")
    50 = @("0.06572", "  -2.57%  ")
    51 = @("1.003", "  -33.37%  ")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $dCell = $ws.Range("D$row")
    $origStyle = $dCell.Style
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]
    $dCell.Style = $origStyle

    $eCell = $ws.Range("E$row")
    $eOrigStyle = $eCell.Style
    $eCell.NumberFormat = "@"
    $eCell.Value = $vals[1]
    $eCell.Style = $eOrigStyle
}
